$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '44.068.87'
$ws.Range("E2").Value = '  -0.96%  '
$ws.Range("D3").Value = '2.244.55'
$ws.Range("E3").Value = '  -1.67%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = "'315.15"
$ws.Range("E5").Value = '  -1.91%  '
$ws.Range("D6").Value = "'99.34"
$ws.Range("E6").Value = '  -6.46%  '
$ws.Range("D7").Value = "'0.574"
$ws.Range("E7").Value = '  -3.00%  '
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("D9").Value = "'0.532"
$ws.Range("E9").Value = '  -7.11%  '
$ws.Range("D10").Value = "'36.16"
$ws.Range("E10").Value = '  -6.89%  '
$ws.Range("D11").Value = "'0.0824"
$ws.Range("E11").Value = '  -2.62%  '
$ws.Range("D12").Value = "'7.36"
$ws.Range("E12").Value = '  -7.11%  '
$ws.Range("E13").Value = '  -2.94%  '
$ws.Range("D14").Value = '2.585.31'
$ws.Range("E14").Value = '  -1.71%  '
$ws.Range("D15").Value = "'0.842"
$ws.Range("E15").Value = '  -5.12%  '
$ws.Range("D16").Value = '2.244.66'
$ws.Range("E16").Value = '  -4.00%  '
$ws.Range("D17").Value = "'13.96"
$ws.Range("E17").Value = '  -4.83%  '
$ws.Range("D18").Value = '43.907.33'
$ws.Range("E18").Value = '  -0.94%  '
$ws.Range("D19").Value = "'13.14"
$ws.Range("E19").Value = '  -7.43%  '
$ws.Range("D20").Value = '0.0₃0978'
$ws.Range("E20").Value = '  -2.84%  '
$ws.Range("D21").Value = "'6.32"
$ws.Range("E21").Value = '  -3.69%  '
$ws.Range("D22").Value = "'65.58"
$ws.Range("E22").Value = '  -1.44%  '
$ws.Range("D23").Value = "'236.51"
$ws.Range("E23").Value = '  -1.29%  '
$ws.Range("D24").Value = "'2.98"
$ws.Range("E24").Value = '  -7.53%  '
$ws.Range("D25").Value = "'2.02"
$ws.Range("E25").Value = '  -8.68%  '
$ws.Range("E26").Value = '  +0.40%  '
$ws.Range("D27").Value = "'10.17"
$ws.Range("E27").Value = '  -0.58%  '
$ws.Range("E28").Value = '  -4.45%  '
$ws.Range("D29").Value = "'36.42"
$ws.Range("E29").Value = '  -5.89%  '
$ws.Range("D30").Value = "'5.98"
$ws.Range("E30").Value = '  -8.97%  '
$ws.Range("D31").Value = "'20.08"
$ws.Range("E31").Value = '  -2.91%  '
$ws.Range("D32").Value = "'155.56"
$ws.Range("E32").Value = '  -4.54%  '
$ws.Range("D33").Value = "'0.0837"
$ws.Range("E33").Value = '  -5.96%  '
$ws.Range("D34").Value = "'3.31"
$ws.Range("E34").Value = '  +2.81%  '
$ws.Range("E35").Value = '  -2.74%  '
$ws.Range("D36").Value = "'1.91"
$ws.Range("E36").Value = '  -8.00%  '
$ws.Range("D37").Value = "'0.108"
$ws.Range("E37").Value = '  -8.64%  '
$ws.Range("D38").Value = "'0.118"
$ws.Range("E38").Value = '  -3.26%  '
$ws.Range("D39").Value = "'15.49"
$ws.Range("E39").Value = '  -0.35%  '
$ws.Range("B40").Value = 'NEARProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D40").Value = "'3.52"
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").Value = "'3.99"
$ws.Range("E41").Value = '  -10.83%  '
$ws.Range("D42").Value = "'0.0307"
$ws.Range("E42").Value = '  -6.92%  '
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").Value = '1.700.54'
$ws.Range("E44").Value = '  -4.33%  '
$ws.Range("D45").Value = "'82.80"
$ws.Range("E45").Value = '  -5.00%  '
$ws.Range("D46").Value = "'0.195"
$ws.Range("E46").Value = '  -6.68%  '
$ws.Range("D47").Value = "'5.18"
$ws.Range("E47").Value = '  -6.14%  '
$ws.Range("D48").Value = "'101.77"
$ws.Range("E48").Value = '  -2.80%  '
$ws.Range("D49").Value = "'71.41"
$ws.Range("E49").Value = '  -4.36%  '
$ws.Range("D50").Value = "'56.29"
$ws.Range("E50").Value = '  -7.10%  '
$ws.Range("D51").Value = "'1.60"
$ws.Range("E51").Value = '  -6.89%  '

# Strip the quote-prefix styling Excel applies for forced-text numeric values
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
